$p = $ppt.ActivePresentation

# Slide 1: update title text
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Intro to Git for Turner Group"

# Slide 7: "How to collaborate using GitHub?" content placeholder text tweaks
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(3).TextFrame.TextRange
$tr7.Runs(2, 1).Text = "Check current branch and other available branches"
$tr7.Runs(6, 1).Text = "Hop to the desired branch"

# Slide 8: "Group Exercise!" content placeholder - insert new numbered bullet
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(3).TextFrame.TextRange
$para5 = $tr8.Paragraphs(5, 1)
[void]$para5.InsertBefore("Make a branch!`r")
